# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 00:23"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5694181
$ws.Range("C4").Value = 38207
$ws.Range("D4").Value = 3053439
$ws.Range("E4").Value = 2464595
$ws.Range("G4").Value = 1073
$ws.Range("H4").Value = 176147

# Row 5: Brasil
$ws.Range("B5").Value = 3456652
$ws.Range("C5").Value = 44780
$ws.Range("D5").Value = 2615254
$ws.Range("E5").Value = 730298
$ws.Range("G5").Value = 1081
$ws.Range("H5").Value = 111100

# Row 27: Canada
$ws.Range("B27").Value = 123422
$ws.Range("C27").Value = 268
$ws.Range("D27").Value = 109777
$ws.Range("E27").Value = 4596
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 9049

# Row 33: Egipto
$ws.Range("B33").Value = 96914
$ws.Range("C33").Value = 161
$ws.Range("D33").Value = 62553
$ws.Range("E33").Value = 29164
$ws.Range("G33").Value = 13
$ws.Range("H33").Value = 5197

# Row 52: Nigeria
$ws.Range("B52").Value = 50488
$ws.Range("C52").Value = 593
$ws.Range("D52").Value = 37304
$ws.Range("E52").Value = 12199
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 985

# Row 61: Uzbekistan
$ws.Range("B61").Value = 37112
$ws.Range("C61").Value = 760
$ws.Range("D61").Value = 32944
$ws.Range("E61").Value = 3920
$ws.Range("G61").Value = 6
$ws.Range("H61").Value = 248

# Row 97: Haiti
$ws.Range("B97").Value = 7949
$ws.Range("C97").Value = 28
$ws.Range("D97").Value = 5337
$ws.Range("E97").Value = 2416

# Row 127: Ruanda
$ws.Range("B127").Value = 2644
$ws.Range("C127").Value = 67
$ws.Range("D127").Value = 1698
$ws.Range("E127").Value = 936

# Row 163: Guyana
$ws.Range("B163").Value = 776
$ws.Range("C163").Value = 39
$ws.Range("D163").Value = 381
$ws.Range("E163").Value = 368
$ws.Range("G163").Value = 2
$ws.Range("H163").Value = 27

# Rows 213/214: swap Islas Malvinas <-> Montserrat (and their stats)
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
